$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160, pushing existing rows 160-285 down to 161-286.
$ws.Rows("160:160").Insert()

# Fill the new row 160 with data. Most columns (A,B,C,E,F,G,H,I,N,O,Q,R) mirror
# what is now row 161 (the original row 160 content), only D,J,K,L,M,P differ.
$ws.Range("A160").Value = 8
$ws.Range("B160").Value = "Terminal La Palmera de La Serena"
$ws.Range("C160").Value = "Coquimbo"
$ws.Range("D160").Value = 44741
$ws.Range("E160").Value = 4
$ws.Range("F160").Value = 100112003
$ws.Range("G160").Value = "Ajo"
$ws.Range("H160").Value = "Chino"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 400
$ws.Range("K160").Value = 18500
$ws.Range("L160").Value = 19000
$ws.Range("M160").Value = 18750
$ws.Range("N160").Value = '$/caja 10 kilos'
$ws.Range("O160").Value = "China"
$ws.Range("P160").Value = 1875
$ws.Range("Q160").Value = 10
$ws.Range("R160").Value = "Hortaliza"
